$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32 (pushes existing row 32..99 down to 33..100)
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new weekly data point
$ws.Cells.Item(32, 1).Value = 4
$ws.Cells.Item(32, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(32, 3).Value = "Los Lagos"
$ws.Cells.Item(32, 4).Value = 44614
$ws.Cells.Item(32, 5).Value = 10
$ws.Cells.Item(32, 6).Value = 100112052
$ws.Cells.Item(32, 7).Value = "Albahaca"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 120
$ws.Cells.Item(32, 11).Value = 5000
$ws.Cells.Item(32, 12).Value = 5000
$ws.Cells.Item(32, 13).Value = 5000
$ws.Cells.Item(32, 14).Value = "$/docena de matas"
$ws.Cells.Item(32, 15).Value = "Región Metropolitana"
$ws.Cells.Item(32, 16).Value = 833
$ws.Cells.Item(32, 17).Value = 6
$ws.Cells.Item(32, 18).Value = "Hortaliza"
